# Update cryptocurrency price/volume data and re-sequence several rows
# per the latest GitHub Actions symbol-list refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cell, $val) {
    # Force the cell to stay a text value (many of these look numeric,
    # e.g. "27.29" or "-1.07%", and Excel would otherwise silently
    # convert them to numbers on assignment).
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "E2" '-1.07%'
Set-TextValue $ws "D3" '27.29'
Set-TextValue $ws "E3" '3.75%'
Set-TextValue $ws "D4" '5.042'
Set-TextValue $ws "E4" '-0.83%'
Set-TextValue $ws "D5" '0.05677'
Set-TextValue $ws "E5" '1.05%'
Set-TextValue $ws "D6" '6.472'
Set-TextValue $ws "E6" '-0.38%'
Set-TextValue $ws "D7" '0.8227'
Set-TextValue $ws "E7" '1.20%'
Set-TextValue $ws "D8" '0.8436'
Set-TextValue $ws "E8" '-0.40%'
Set-TextValue $ws "B9" 'WazirX'
Set-TextValue $ws "C9" 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws "D9" '0.1326'
Set-TextValue $ws "E9" '-1.49%'
Set-TextValue $ws "B10" 'MandalaExchangeToken'
Set-TextValue $ws "C10" 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws "D10" '0.06914'
Set-TextValue $ws "E10" '-0.75%'
Set-TextValue $ws "B11" 'LiechtensteinCryptoassetsExchange'
Set-TextValue $ws "C11" 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws "D11" '0.03167'
Set-TextValue $ws "E11" '0.09%'
Set-TextValue $ws "D12" '0.02884'
Set-TextValue $ws "E12" '2.20%'
Set-TextValue $ws "E13" '-0.15%'
Set-TextValue $ws "D14" '0.001509'
Set-TextValue $ws "E14" '-0.84%'
Set-TextValue $ws "D15" '0.04137'
Set-TextValue $ws "E15" '-11.59%'
Set-TextValue $ws "B16" 'One'
Set-TextValue $ws "C16" 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue $ws "D16" '0.0006009'
Set-TextValue $ws "E16" '0.45%'
Set-TextValue $ws "B17" 'TigerCash'
Set-TextValue $ws "C17" 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws "D17" '0.006180'
Set-TextValue $ws "E17" '-0.94%'
Set-TextValue $ws "B18" 'LEO'
Set-TextValue $ws "C18" 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws "D18" '3.513'
Set-TextValue $ws "E18" '-1.67%'
Set-TextValue $ws "B19" 'GateToken'
Set-TextValue $ws "C19" 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws "D19" '2.999'
Set-TextValue $ws "E19" '-1.85%'
Set-TextValue $ws "B20" 'BTSEToken'
Set-TextValue $ws "C20" 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws "D20" '2.225'
Set-TextValue $ws "E20" '5.06%'
Set-TextValue $ws "B21" 'BitpandaEcosystemToken'
Set-TextValue $ws "C21" 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue $ws "D21" '0.3113'
Set-TextValue $ws "E21" '-2.14%'
Set-TextValue $ws "E22" '-5.00%'
Set-TextValue $ws "D23" '3.560'
Set-TextValue $ws "E23" '-5.11%'
Set-TextValue $ws "D25" '0.001219'
Set-TextValue $ws "E25" '-2.37%'
Set-TextValue $ws "D26" '0.003869'
Set-TextValue $ws "E26" '-16.29%'
Set-TextValue $ws "E27" '2.04%'
Set-TextValue $ws "E28" '-25.80%'
Set-TextValue $ws "D40" '0.03664'
Set-TextValue $ws "E40" '-0.01%'
Set-TextValue $ws "B41" 'BKEXToken'
Set-TextValue $ws "C41" 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue $ws "D41" '0.1054'
Set-TextValue $ws "E41" '-0.48%'
Set-TextValue $ws "B42" 'KickToken'
Set-TextValue $ws "C42" 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue $ws "D42" '0.006050'
Set-TextValue $ws "E42" '-1.19%'
Set-TextValue $ws "D43" '0.002278'
Set-TextValue $ws "E43" '-12.35%'
Set-TextValue $ws "D44" '0.009610'
Set-TextValue $ws "E44" '10.49%'
Set-TextValue $ws "D45" '0.00005310'
Set-TextValue $ws "E45" '0.29%'
Set-TextValue $ws "E46" '-0.05%'
Set-TextValue $ws "E47" '-15.45%'
Set-TextValue $ws "D48" '0.002561'
Set-TextValue $ws "E48" '23.92%'
Set-TextValue $ws "E49" '-0.05%'
Set-TextValue $ws "E50" '-0.05%'